$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits after
#    "End Planet." (it will be re-created at the new location below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the paragraph containing "Joining of results 2 & 3" and
#    rewrite it as four separate runs:
#       "Joining of results " | "1" | " & " | "2"
#    followed by a bookmarkStart/bookmarkEnd pair named "_GoBack".
# ------------------------------------------------------------------
$target = $d.Range(0, 0)
$found = $target.Find.Execute("Joining of results 2 & 3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Joining of results 2 & 3' paragraph"
}

# A throwaway bookmark placed before any edits keeps the engine from
# re-coalescing the runs we are about to insert separately.
$guard = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("zzz_run_guard", $guard)

$r = $d.Range($target.Start, $target.End)
$r.Text = ""
$r.InsertAfter("Joining of results ")
$r.Collapse(0)
$r.InsertAfter("1")
$r.Collapse(0)
$r.InsertAfter(" & ")
$r.Collapse(0)
$r.InsertAfter("2")
$r.Collapse(0)

# Insert a throw-away trailing character so the new "_GoBack" bookmark
# can be placed *between* two runs rather than snapping to the edge of
# the paragraph (which the engine treats specially).
$r.InsertAfter("X")
$r.Collapse(0)

$bmRange = $d.Range($r.Start - 1, $r.Start - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the throw-away trailing character again.
$dummy = $d.Range($r.Start - 1, $r.Start)
$dummy.Text = ""

# Remove the guard bookmark used only to stop run coalescing.
if ($d.Bookmarks.Exists("zzz_run_guard")) {
    $d.Bookmarks("zzz_run_guard").Delete()
}
